# Generate Report for Handoff
#
# The "5881c840-64b8-412f-ba32-66daa693f7d6" file moves from "In Translation"
# to "Ready for handoff" across all three sheets, and its handoff
# date/datetime stamps are refreshed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-13-14 02:13:26"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-14 02:13:24"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-14 02:13:26"
